# Define build_type 2 as office building. (#2)
#
# The non-residential specific-demand lookup table (sheet "Tabelle1") had a
# row inserted conceptually in front of the data: every existing
# "type_id" in column A (rows 2-45) is bumped up by one so that, elsewhere
# in the codebase, build_type == 2 can be (re-)defined to mean an office
# building. Columns B-F (name / demands / SLP types) keep referring to the
# same row they always did - only the numeric id in column A shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row r (2..45) previously held type_id = r - 1; it now holds type_id = r.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 1).Value = $r
}

# Reflect where the workbook was left scrolled/selected when last saved.
$ws.Range("B30").Select()
